# Weekly update: a new observation is inserted at row 189 (shifting the
# existing rows 189-275 down to 190-276), matching the "Fruta / hortaliza,
# semanal" refresh pattern used across these consolidated sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by one row, starting at row 189.
$ws.Rows("189:189").Insert()

# Populate the newly inserted row with this week's observation.
$ws.Range("A189").Value = 6
$ws.Range("B189").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C189").Value = "Metropolitana"
$ws.Range("D189").Value = 44917
$ws.Range("E189").Value = 13
$ws.Range("F189").Value = 100112029
$ws.Range("G189").Value = "Orégano"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 46
$ws.Range("K189").Value = 16000
$ws.Range("L189").Value = 17000
$ws.Range("M189").Value = 16457
$ws.Range("N189").Value = "$/docena de atados"
$ws.Range("O189").Value = "Región Metropolitana"
$ws.Range("P189").Value = 5486
$ws.Range("Q189").Value = 3
$ws.Range("R189").Value = "Hortaliza"
